$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Running all the test cases: set Runmode (column D) to "Y" for every test case row
$range = $ws.Range("D2:D25")
$range.Select()
$range.Value = "Y"

# Update Results (column E) to reflect the outcome of actually running the
# previously-skipped test cases
$ws.Range("E11").Value = "PASS"
$ws.Range("E13").Value = "SKIP"
